$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated calibration values after removing the < USD 5 price noise point
# from the extrapolation calibration. Columns: D (ABSM1_RN), E (M1_RN),
# F (CM2_RN), G (CMN3_RN), H (CMN4_RN).

$updates = @{
    "D5"  = 118813.4402538303
    "E5"  = 0.0009396634092559489
    "F5"  = 0.2562687968355845
    "G5"  = -3.166690184265289
    "H5"  = 30.69701035915155

    "D6"  = 119426.2414862402
    "E6"  = -0.003252274816201206
    "F6"  = 0.2306427786821414
    "G6"  = -1.37460089863363
    "H6"  = 11.09111437016731

    "D7"  = 120778.5835379928
    "E7"  = -0.007682084482526705
    "F7"  = 0.2397235403484029
    "G7"  = -0.8618971446192862
    "H7"  = 6.029805657320228

    "D8"  = 121158.5365706338
    "E8"  = -0.02934473453374065
    "F8"  = 0.2089456869718486
    "G8"  = -0.8666898587295436
    "H8"  = 6.408555039204717

    "D9"  = 122871.8380957177
    "E9"  = -0.06077589391522076
    "F9"  = 0.3321942907472831
    "G9"  = -1.564542624292799
    "H9"  = 10.24899334351268

    "D10" = 124274.877566475
    "E10" = -0.09964543623477756
    "F10" = 0.4329951691417186
    "G10" = -1.904215885590203
    "H10" = 9.830522351798576

    "D11" = 126301.6564948566
    "E11" = -0.1740458614175637
    "F11" = 0.7439092569507443
    "G11" = -2.56451977851868
    "H11" = 12.40393845884312

    "D16" = 116642.6763131885
    "E16" = 0.1254973488098925
    "F16" = 0.131530660456945
    "G16" = -0.3791559876162946
    "H16" = 5.590775207335812

    "D17" = 116591.6017221299
    "E17" = 0.102464528396485
    "F17" = 0.1388459732910115
    "G17" = -0.6286019549745252
    "H17" = 7.545776870995195

    "D19" = 117420.7945423491
    "E19" = 0.05151742074064533
    "F19" = 0.1456571385932011
    "G19" = -0.5585721077709873
    "H19" = 6.340494539439065

    "D20" = 117600.2326832831
    "E20" = 0.04547191680145291
    "F20" = 0.1462089083649787
    "G20" = -0.5361635184772116
    "H20" = 6.771562382336227
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
